$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date updated
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was blank, now "Alvearie Team"
$ws.Range("B9").Value = "Alvearie Team"

# Row 10: Contact / No display for ContactDetail -> Jurisdiction / United States of America
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 was a duplicate Contact row; remove it entirely, shifting rows 12-15 up
$ws.Rows.Item(11).Delete()
